$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 7594
}

for ($r = 14; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 7586
}

for ($r = 27; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7569
}
